$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 4,23
$data[0,0] = 0.0045489006823351
$data[0,1] = 0.00985595147839272
$data[0,2] = 0.0128885519332828
$data[0,3] = 0.943896891584534
$data[0,4] = 0.0151630022744503
$data[0,5] = 0.980288097043215
$data[0,6] = 0.967399545109932
$data[0,7] = 0.0106141015921152
$data[0,8] = 0.00606520090978014
$data[0,9] = 0.00151630022744503
$data[0,10] = 0.99090219863533
$data[0,11] = 0
$data[0,12] = 0
$data[0,13] = 0
$data[0,14] = 0
$data[0,15] = 0.000758150113722517
$data[0,16] = 0.000758150113722517
$data[0,17] = 0.0363912054586808
$data[0,18] = 0.0295678544351782
$data[0,19] = 0.00227445034116755
$data[0,20] = 0.0151630022744503
$data[0,21] = 0.0166793025018954
$data[0,22] = 0.0045489006823351
$data[1,0] = 0.0181956027293404
$data[1,1] = 0.962850644427597
$data[1,2] = 0.00682335102350265
$data[1,3] = 0.00682335102350265
$data[1,4] = 0.970432145564822
$data[1,5] = 0.0174374526156179
$data[1,6] = 0.00227445034116755
$data[1,7] = 0
$data[1,8] = 0.00227445034116755
$data[1,9] = 0
$data[1,10] = 0
$data[1,11] = 0.998483699772555
$data[1,12] = 0.00379075056861259
$data[1,13] = 0.0045489006823351
$data[1,14] = 0.00227445034116755
$data[1,15] = 0.99696739954511
$data[1,16] = 0.989385898407885
$data[1,17] = 0
$data[1,18] = 0.954510993176649
$data[1,19] = 0.020470053070508
$data[1,20] = 0.00151630022744503
$data[1,21] = 0
$data[1,22] = 0
$data[2,0] = 0.187263078089462
$data[2,1] = 0.00303260045489007
$data[2,2] = 0.0128885519332828
$data[2,3] = 0.043972706595906
$data[2,4] = 0.00379075056861259
$data[2,5] = 0.00151630022744503
$data[2,6] = 0.0250189537528431
$data[2,7] = 0.988627748294162
$data[2,8] = 0.991660348749052
$data[2,9] = 0.995451099317665
$data[2,10] = 0.00530705079605762
$data[2,11] = 0
$data[2,12] = 0.000758150113722517
$data[2,13] = 0
$data[2,14] = 0.000758150113722517
$data[2,15] = 0
$data[2,16] = 0.000758150113722517
$data[2,17] = 0.963608794541319
$data[2,18] = 0.00379075056861259
$data[2,19] = 0.000758150113722517
$data[2,20] = 0.974981046247157
$data[2,21] = 0.982562547384382
$data[2,22] = 0.991660348749052
$data[3,0] = 0.789992418498863
$data[3,1] = 0.0242608036391205
$data[3,2] = 0.967399545109932
$data[3,3] = 0.00530705079605762
$data[3,4] = 0.0106141015921152
$data[3,5] = 0.000758150113722517
$data[3,6] = 0.00530705079605762
$data[3,7] = 0
$data[3,8] = 0
$data[3,9] = 0.00303260045489007
$data[3,10] = 0.00379075056861259
$data[3,11] = 0.000758150113722517
$data[3,12] = 0.995451099317665
$data[3,13] = 0.995451099317665
$data[3,14] = 0.99696739954511
$data[3,15] = 0.00227445034116755
$data[3,16] = 0.0090978013646702
$data[3,17] = 0
$data[3,18] = 0.0121304018195603
$data[3,19] = 0.976497346474602
$data[3,20] = 0.00758150113722517
$data[3,21] = 0
$data[3,22] = 0.00303260045489007

$ws.Range("B2:X5").Value = $data
